$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GitHub activity numbers
$ws.Range("C8").Value = 5
$ws.Range("C9").Value = 26

# Basic options scores that were previously blank
$ws.Range("C14").Value = 0
$ws.Range("C19").Value = 10
$ws.Range("C20").Value = 5
$ws.Range("C21").Value = 5
$ws.Range("C22").Value = 5
$ws.Range("C32").Value = 5

# Move the selection / scroll position like the author's final view
$ws.Range("C11").Select()
